# feat(data): llm qwen update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update pricing (input/output token price) for existing Qwen rows 51-61
#    (columns G = input token price, H = output token price)
# ---------------------------------------------------------------------------
$priceUpdates = @{
    51 = @(0.14,               0.28)
    52 = @(0.28,                0.56)
    53 = @(0.49,                0.98)
    54 = @(0.7,                 1.4)
    55 = @(0.98,                1.96)
    56 = @(0.28,                0.84)
    57 = @(0.56,                1.68)
    58 = @(16.8,                16.8)
    59 = @(5.6,                 16.8)
    60 = @(5.6,                 16.8)
    61 = @(5.6,                 16.8)
}

foreach ($row in $priceUpdates.Keys) {
    $vals = $priceUpdates[$row]
    $ws.Cells.Item($row, 7).Value2 = $vals[0]
    $ws.Cells.Item($row, 8).Value2 = $vals[1]
}

# ---------------------------------------------------------------------------
# 2) Insert two new rows (62 and 63) for the new "Qwen Long" and "Qwen VL"
#    models, pushing all the following rows down by two.
# ---------------------------------------------------------------------------
$ws.Rows("62:63").Insert()

# New row 62: Qwen Long
$ws.Cells.Item(62, 1).Value2 = "Qwen Long"
$ws.Cells.Item(62, 2).Value2 = 45433
$ws.Cells.Item(62, 3).Value2 = "Aliyun"
$ws.Cells.Item(62, 4).Value2 = "110*"
$ws.Cells.Item(62, 5).Value2 = 10000
$ws.Cells.Item(62, 6).Value2 = 2
$ws.Cells.Item(62, 7).Value2 = 0.07
$ws.Cells.Item(62, 8).Value2 = 0.28
$ws.Cells.Item(62, 11).Value2 = "Proprietary"
$ws.Cells.Item(62, 12).Value2 = $true
$ws.Cells.Item(62, 13).Value2 = $true

# New row 63: Qwen VL
$ws.Cells.Item(63, 1).Value2 = "Qwen VL"
$ws.Cells.Item(63, 2).Value2 = 45217
$ws.Cells.Item(63, 3).Value2 = "Aliyun"
$ws.Cells.Item(63, 4).Value2 = 7
$ws.Cells.Item(63, 5).Value2 = 6
$ws.Cells.Item(63, 6).Value2 = 2
$ws.Cells.Item(63, 7).Value2 = 0
$ws.Cells.Item(63, 8).Value2 = 0
$ws.Cells.Item(63, 9).Value2 = 0
$ws.Cells.Item(63, 11).Value2 = "tongyi-qianwen"
$ws.Cells.Item(63, 12).Value2 = $true
$ws.Cells.Item(63, 14).Value2 = $true
$ws.Cells.Item(63, 16).Value2 = 58.2
$ws.Cells.Item(63, 18).Value2 = 31.6
$ws.Cells.Item(63, 19).Value2 = 11.6

# ---------------------------------------------------------------------------
# 3) "Qwen VL Max" (now row 65) is no longer the first of its series, since
#    "Qwen VL" now holds that spot - clear the "series first" flag (col L).
# ---------------------------------------------------------------------------
$ws.Cells.Item(65, 12).ClearContents()

# ---------------------------------------------------------------------------
# 4) Update the frozen-pane view state to match the author's final selection.
# ---------------------------------------------------------------------------
$ws.Range("H56").Select()
$excel.ActiveWindow.ScrollRow = 42
